$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update story title and link (source-name swap + new tinyurl)
$ws.Range("B5").Value = "Three-year drug investigation in Chautauqua Co. ends in 2 lb. meth bust - InformNNY.com"
$ws.Range("E5").Value = "https://tinyurl.com/22yfkywo"

# Row 15 (old "WNY News Now / Federal Grand Jury Indicts..." entry) is removed entirely;
# deleting it shifts rows 16-20 up to 15-19, matching the rest of the diff.
$ws.Rows("15:15").Delete()

# After the shift, what was row 18 is now row 17; update it to the "Rolison Sentenced..." story.
$ws.Range("B17").Value = "Rolison Sentenced For Having Meth At Jail - Post Journal"
$ws.Range("D17").Value = 45748
$ws.Range("E17").Value = "https://tinyurl.com/2amvc7b2"

# What was row 19 is now row 18; update it to the "Police raid at Dunkirk home..." story.
$ws.Range("B18").Value = "Police raid at Dunkirk home leads to meth-dealing charges against 2 - The Star Press"
$ws.Range("D18").Value = 45748
$ws.Range("E18").Value = "https://tinyurl.com/2xjgrbsq"
